$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3.494130772323899
$ws.Range("D2").Value = 4.672245578028645
$ws.Range("E2").Value = 10.72125653781629
$ws.Range("F2").Value = 24.9379913592078
$ws.Range("G2").Value = 30.20748193470057
$ws.Range("H2").Value = 14.37780188253867
$ws.Range("I2").Value = 20.6417696799914
$ws.Range("K2").Value = 13.54477112006325
$ws.Range("M2").Value = 15.40947398701308

$ws.Range("C3").Value = 3.414013567003487
$ws.Range("D3").Value = 4.666857086176097
$ws.Range("E3").Value = 10.65758884208388
$ws.Range("F3").Value = 24.92639173636222
$ws.Range("G3").Value = 30.13935031644152
$ws.Range("H3").Value = 14.43623324354205
$ws.Range("I3").Value = 20.6769484384216
$ws.Range("K3").Value = 12.9738543449897
$ws.Range("M3").Value = 15.0921539315271

$ws.Range("C4").Value = 3.363422214283655
$ws.Range("D4").Value = 4.663833046233505
$ws.Range("E4").Value = 10.62299286192739
$ws.Range("F4").Value = 24.92990634986318
$ws.Range("G4").Value = 30.11413693792126
$ws.Range("H4").Value = 14.47611594325838
$ws.Range("I4").Value = 20.70634394374765
$ws.Range("K4").Value = 12.611451105029
$ws.Range("M4").Value = 14.89735900020397

$ws.Range("C5").Value = 3.342472347500905
$ws.Range("D5").Value = 4.662673076176831
$ws.Range("E5").Value = 10.61003428385933
$ws.Range("F5").Value = 24.93400347181082
$ws.Range("G5").Value = 30.10803028906526
$ws.Range("H5").Value = 14.49337088510514
$ws.Range("I5").Value = 20.72027159302756
$ws.Range("K5").Value = 12.46099294812715
$ws.Range("M5").Value = 14.81810142928583

$ws.Range("C6").Value = 3.338974054252307
$ws.Range("D6").Value = 4.662484857284311
$ws.Range("E6").Value = 10.60795156179269
$ws.Range("F6").Value = 24.93484440075881
$ws.Range("G6").Value = 30.10726754672011
$ws.Range("H6").Value = 14.49629645588348
$ws.Range("I6").Value = 20.72270160660072
$ws.Range("K6").Value = 12.43584836798674
$ws.Range("M6").Value = 14.80495153127732

$ws.Range("C7").Value = 3.363141002273546
$ws.Range("D7").Value = 4.663817108489829
$ws.Range("E7").Value = 10.62281347398639
$ws.Range("F7").Value = 24.92995083043905
$ws.Range("G7").Value = 30.11403772596294
$ws.Range("H7").Value = 14.47634459718655
$ws.Range("I7").Value = 20.70652390275199
$ws.Range("K7").Value = 12.60943292164753
$ws.Range("M7").Value = 14.89628945659395

$ws.Range("C8").Value = 3.466807378733059
$ws.Range("D8").Value = 4.670329059141368
$ws.Range("E8").Value = 10.69837683984362
$ws.Range("F8").Value = 24.93177916754029
$ws.Range("G8").Value = 30.18053637467053
$ws.Range("H8").Value = 14.39711477711846
$ws.Range("I8").Value = 20.65227478487001
$ws.Range("K8").Value = 13.35049562782715
$ws.Range("M8").Value = 15.30012038604344

$ws.Range("C9").Value = 3.658249548740771
$ws.Range("D9").Value = 4.685321552770036
$ws.Range("E9").Value = 10.88172066214005
$ws.Range("F9").Value = 25.02005248344666
$ws.Range("G9").Value = 30.44292801462057
$ws.Range("H9").Value = 14.27374006927208
$ws.Range("I9").Value = 20.60820250904826
$ws.Range("K9").Value = 14.7018683480434
$ws.Range("M9").Value = 16.08754730674318

$ws.Range("C10").Value = 3.790763595016394
$ws.Range("D10").Value = 4.697645382844795
$ws.Range("E10").Value = 11.03701430325652
$ws.Range("F10").Value = 25.13671089341696
$ws.Range("G10").Value = 30.71578667293488
$ws.Range("H10").Value = 14.20288579795416
$ws.Range("I10").Value = 20.61434042141905
$ws.Range("K10").Value = 15.62386629588693
$ws.Range("M10").Value = 16.65715946334989

$ws.Range("C11").Value = 3.849109995679255
$ws.Range("D11").Value = 4.703526130508123
$ws.Range("E11").Value = 11.11190766603669
$ws.Range("F11").Value = 25.20099782047386
$ws.Range("G11").Value = 30.85705832065243
$ws.Range("H11").Value = 14.17500805795168
$ws.Range("I11").Value = 20.62557969162646
$ws.Range("K11").Value = 16.02648064539623
$ws.Range("M11").Value = 16.91314030764267

$ws.Range("C12").Value = 3.870913283447285
$ws.Range("D12").Value = 4.705791633936835
$ws.Range("E12").Value = 11.14085723761596
$ws.Range("F12").Value = 25.22694663086422
$ws.Range("G12").Value = 30.91298759236392
$ws.Range("H12").Value = 14.16508225367095
$ws.Range("I12").Value = 20.63105471061141
$ws.Range("K12").Value = 16.17642653342772
$ws.Range("M12").Value = 17.00952818812696

$ws.Range("C13").Value = 3.866230702814281
$ws.Range("D13").Value = 4.705302015103498
$ws.Range("E13").Value = 11.13459659709513
$ws.Range("F13").Value = 25.22128688912438
$ws.Range("G13").Value = 30.90083462666665
$ws.Range("H13").Value = 14.16719181995681
$ws.Range("I13").Value = 20.62982130342211
$ws.Range("K13").Value = 16.14424628751716
$ws.Range("M13").Value = 16.98879509288406

$ws.Range("C14").Value = 3.850909672392669
$ws.Range("D14").Value = 4.703711746527717
$ws.Range("E14").Value = 11.11427769035116
$ws.Range("F14").Value = 25.20310054361108
$ws.Range("G14").Value = 30.86161111375136
$ws.Range("H14").Value = 14.17417878533518
$ws.Range("I14").Value = 20.62600567645837
$ws.Range("K14").Value = 16.03886770770874
$ws.Range("M14").Value = 16.921081673846

$ws.Range("C15").Value = 3.841486805820619
$ws.Range("D15").Value = 4.702742660350692
$ws.Range("E15").Value = 11.10190780861194
$ws.Range("F15").Value = 25.19216955878241
$ws.Range("G15").Value = 30.83790129550517
$ws.Range("H15").Value = 14.17854080271289
$ws.Range("I15").Value = 20.62382733719797
$ws.Range("K15").Value = 15.97399003531795
$ws.Range("M15").Value = 16.87953135262467

$ws.Range("C16").Value = 3.786910449309773
$ws.Range("D16").Value = 4.697266511219551
$ws.Range("E16").Value = 11.03220351437295
$ws.Range("F16").Value = 25.13273468997657
$ws.Range("G16").Value = 30.7068966296605
$ws.Range("H16").Value = 14.20479570202869
$ws.Range("I16").Value = 20.61377628776415
$ws.Range("K16").Value = 15.59720751947583
$ws.Range("M16").Value = 16.64035908737608

$ws.Range("C17").Value = 3.752924448382838
$ws.Range("D17").Value = 4.69397676384213
$ws.Range("E17").Value = 10.99051547097069
$ws.Range("F17").Value = 25.09914169638948
$ws.Range("G17").Value = 30.63090014632042
$ws.Range("H17").Value = 14.22202088474932
$ws.Range("I17").Value = 20.60977734908554
$ws.Range("K17").Value = 15.36168383755475
$ws.Range("M17").Value = 16.49276122148814

$ws.Range("C18").Value = 3.733195181749527
$ws.Range("D18").Value = 4.692110486237582
$ws.Range("E18").Value = 10.96693902355396
$ws.Range("F18").Value = 25.08087654997868
$ws.Range("G18").Value = 30.5888053756257
$ws.Range("H18").Value = 14.23233791869568
$ws.Range("I18").Value = 20.60827219319939
$ws.Range("K18").Value = 15.22463972065647
$ws.Range("M18").Value = 16.40757674683019

$ws.Range("C19").Value = 3.726484446543093
$ws.Range("D19").Value = 4.691483071792785
$ws.Range("E19").Value = 10.95902602942254
$ws.Range("F19").Value = 25.0748739686057
$ws.Range("G19").Value = 30.57483130992561
$ws.Range("H19").Value = 14.23590128479162
$ws.Range("I19").Value = 20.6078989463075
$ws.Range("K19").Value = 15.1779713228268
$ws.Range("M19").Value = 16.37868792016904

$ws.Range("C20").Value = 3.756561189892278
$ws.Range("D20").Value = 4.694324289193313
$ws.Range("E20").Value = 10.99491184383591
$ws.Range("F20").Value = 25.10260841884933
$ws.Range("G20").Value = 30.63882302646252
$ws.Range("H20").Value = 14.22014481148872
$ws.Range("I20").Value = 20.61012073506751
$ws.Range("K20").Value = 15.38691972467665
$ws.Range("M20").Value = 16.50850399710716

$ws.Range("C21").Value = 3.855417832284082
$ws.Range("D21").Value = 4.704177806716121
$ws.Range("E21").Value = 11.12023004156952
$ws.Range("F21").Value = 25.2083988458383
$ws.Range("G21").Value = 30.8730662842467
$ws.Range("H21").Value = 14.17210938484039
$ws.Range("I21").Value = 20.62709331062729
$ws.Range("K21").Value = 16.06988892970691
$ws.Range("M21").Value = 16.94098630808078

$ws.Range("C22").Value = 3.918323612336096
$ws.Range("D22").Value = 4.710842175823815
$ws.Range("E22").Value = 11.20555505121494
$ws.Range("F22").Value = 25.28688640956518
$ws.Range("G22").Value = 31.04031748941695
$ws.Range("H22").Value = 14.14439487304274
$ws.Range("I22").Value = 20.64529167954316
$ws.Range("K22").Value = 16.50155151267905
$ws.Range("M22").Value = 17.22041303857533

$ws.Range("C23").Value = 3.884909422345984
$ws.Range("D23").Value = 4.707265028516093
$ws.Range("E23").Value = 11.15971001555828
$ws.Range("F23").Value = 25.24414449143748
$ws.Range("G23").Value = 30.94976949927214
$ws.Range("H23").Value = 14.15884842908872
$ws.Range("I23").Value = 20.63492769608668
$ws.Range("K23").Value = 16.27253861218459
$ws.Range("M23").Value = 17.07160278879966

$ws.Range("C24").Value = 3.754917610682182
$ws.Range("D24").Value = 4.694167094894748
$ws.Range("E24").Value = 10.99292302547614
$ws.Range("F24").Value = 25.10103784888208
$ws.Range("G24").Value = 30.63523611665418
$ws.Range("H24").Value = 14.22099169483775
$ws.Range("I24").Value = 20.60996301772386
$ws.Range("K24").Value = 15.37551567577958
$ws.Range("M24").Value = 16.50138770414811

$ws.Range("C25").Value = 3.607826832659305
$ws.Range("D25").Value = 4.681032473701799
$ws.Range("E25").Value = 10.82843537574903
$ws.Range("F25").Value = 24.9870744235303
$ws.Range("G25").Value = 30.35782239858454
$ws.Range("H25").Value = 14.30366115897566
$ws.Range("I25").Value = 20.61339390205103
$ws.Range("K25").Value = 14.34817812173327
$ws.Range("M25").Value = 15.87566230597939
